$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $old)
    }
}

# Title
Replace-Text "Unraveling the Enigma of Dark Matter" "Mathematics: The Art of Patterns Unveiled"

# Byline name
Replace-Text " Amelia Rodriguez" " Camille Robinson"

# Email local-part and domain
Replace-Text "ameliarodriguez@astronews" "camillerobinson@highschool"
$foundCom = $d.Content.Find.Execute("com", $true, $true, $false, $false, $false, $true, 1, $false, "edu", 2)
if (-not $foundCom) {
    Write-Output "NOT FOUND: com (whole word)"
}

# Paragraph 1 of body (first block of three sentences collapses into two)
Replace-Text "For centuries, astronomers and physicists have grappled with the mystery of dark matter, an invisible substance that exerts a gravitational pull on visible matter, yet remains elusive to direct observation" "Mathematics, akin to a symphony of numbers and symbols, invites us on an intellectual adventure, where beauty and logic harmoniously intertwine"

Replace-Text " Its existence has been inferred through various astrophysical phenomena, such as the rotation curves of galaxies, the dynamics of galaxy clusters, and gravitational lensing. Dark matter is believed to constitute approximately 85% of the total mass of the universe, yet its composition and properties have remained enigmatic. The quest to unravel this cosmic puzzle has propelled advancements in observational techniques, theoretical models, and experimental setups, leading to a deeper understanding of the cosmos" " From the ancient Babylonians and Egyptians who marveled at the patterns of the stars to the modern-day mathematicians delving into the enigmatic world of quantum physics, the allure of mathematics has captivated minds across time and civilizations"

# Paragraph 1, second block
Replace-Text "With its enigmatic nature, dark matter presents a challenge to our current understanding of physics" "Unveiling the hidden complexities of the cosmos, mathematics unveils patterns and symmetries that govern the world around us"

Replace-Text " It challenges the predictions of Newtonian gravity and forces scientists to explore alternative theories such as modified gravity or extra dimensions" " The intricate dance of celestial bodies, the rhythmic ebb and flow of tides, the mesmerizing spirals found in nature--all bear witness to the profound elegance of mathematical principles"

Replace-Text " Unraveling the mystery of dark matter promises to shed light on the fundamental nature of gravity, the evolution of galaxies and cosmic structures, and the ultimate fate of the universe. From the depths of space to the confines of underground laboratories, the quest for answers continues, pushing the boundaries of human knowledge and offering tantalizing glimpses into the unseen forces that shape our universe" " As we embark on this mathematical odyssey, we not only unravel the mysteries of the universe but also cultivate critical thinking skills, problem-solving abilities, and a profound appreciation for the universe"

# Paragraph 1, third block
Replace-Text "As scientists delve deeper into the cosmos, they encounter a landscape of cosmic mysteries" "Mathematics, like a master artist, paints the canvas of our world with patterns, harmony, and order"

Replace-Text " Dark matter stands as one of the most confounding enigmas, challenging our understanding of the universe" " From the golden ratio in art and architecture to the Fibonacci sequence in biology, mathematics manifests itself in diverse and fascinating ways"

Replace-Text " Its presence, inferred through gravitational effects, hints at an underlying reality that remains obscured from direct observation. Unveiling the secrets of dark matter promises to rewrite our textbooks, redefine our comprehension of gravity, and open new vistas of knowledge in physics and cosmology. This elusive cosmic entity holds the key to unlocking profound mysteries about the nature of matter, the evolution of the cosmos, and the ultimate fate of our universe" " Through this exploration, we embark on a journey of discovery, unraveling the intricate connections between seemingly disparate fields, revealing the interconnectedness of the universe, and fostering a sense of unity and wonder"

# Summary paragraph
Replace-Text "The enigma of dark matter has captivated the scientific community for decades" "Mathematics, an exquisite tapestry of patterns and logic, unveils the hidden order of the universe, inviting us to explore its intricate beauty and underlying harmonies"

Replace-Text " Its existence, inferred through gravitational effects, has challenged our understanding of physics and pushed the boundaries of human knowledge" " As we embark on this mathematical adventure, we not only cultivate critical thinking and problem-solving skills but also deepen our appreciation for the elegance and interconnectedness of the world around us"

Replace-Text " As we delve deeper into the cosmos, unraveling the secrets of dark matter promises to rewrite our textbooks, redefine our comprehension of gravity, and unveil the mysteries of the unseen forces that shape our universe" " This exploration of patterns, symmetries, and relationships equips us with the tools to decipher the enigmatic world of mathematics, unlocking its secrets and unraveling its captivating mysteries"

# Append a new empty paragraph at the very end of the document body
$d.Paragraphs.Add() | Out-Null
